# Daily attendance processing - 2025-10-28 20:24:21
# Normalize the "Recorded By" (column G) lists: entries recorded by the
# automated "System" account should be listed first, followed by the
# remaining recorders in alphabetical order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value()

    if ($null -eq $val) { continue }
    if ($val -eq "") { continue }
    if ($val.IndexOf(",") -lt 0) { continue }

    $parts = $val.Split(",")

    $systemPart = @()
    $otherPart = @()
    foreach ($p in $parts) {
        $trimmed = $p.Trim()
        if ($trimmed.ToLower() -eq "system") {
            $systemPart += $trimmed
        } else {
            $otherPart += $trimmed
        }
    }

    $otherSorted = $otherPart | Sort-Object
    $newParts = $systemPart + $otherSorted
    $newVal = $newParts -join ", "

    if ($newVal -ne $val) {
        $ws.Cells.Item($r, 7).Value = $newVal
    }
}
